$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "Way too many systems but fortunately most systems same as in previous vessels."; New = "Design: Way too many systems but fortunately most systems same as in previous vessels." },
    @{ Old = "More-less straight forward. Only Panama arrangement not fulfilled due to aft ship design."; New = "Design: More-less straight forward. Only Panama arrangement not fulfilled due to aft ship design." },
    @{ Old = "I consider material handling quite smooth. Mostly because it was already third vessel."; New = "Design: I consider material handling quite smooth. Mostly because it was already third vessel." },
    @{ Old = "Block manufacturing timetable and detail design areas could have been more in line to give extra time for design."; New = "Design: Block manufacturing timetable and detail design areas could have been more in line to give extra time for design." },
    @{ Old = "In my opinion communication was smooth between different parties."; New = "Design: In my opinion communication was smooth between different parties." },
    @{ Old = "We got us well employed by changing the hull structure, suppliers and sub-contractors in third vessel."; New = "Design: We got us well employed by changing the hull structure, suppliers and sub-contractors in third vessel." }
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
